$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at 366-368, pushing the existing rows 366-388 down to 369-391.
$ws.Range("A366:R368").EntireRow.Insert()

# Row 366: new "Limache" 1a (cosecha) entry, $/caja 10 kilos
$ws.Range("A366").Value = 3
$ws.Range("B366").Value = "Femacal de La Calera"
$ws.Range("C366").Value = "Coquimbo"
$ws.Range("D366").Value = 44610
$ws.Range("E366").Value = 5
$ws.Range("F366").Value = 100112003
$ws.Range("G366").Value = "Ajo"
$ws.Range("H366").Value = "Chino"
$ws.Range("I366").Value = "1a (cosecha)"
$ws.Range("J366").Value = 70
$ws.Range("K366").Value = 16000
$ws.Range("L366").Value = 16500
$ws.Range("M366").Value = 16250
$ws.Range("N366").Value = "`$/caja 10 kilos"
$ws.Range("O366").Value = "Limache"
$ws.Range("P366").Value = 1625
$ws.Range("Q366").Value = 10
$ws.Range("R366").Value = "Hortaliza"

# Row 367: new "Limache" 1a (cosecha) entry, $/trenza 50 unidades
$ws.Range("A367").Value = 3
$ws.Range("B367").Value = "Femacal de La Calera"
$ws.Range("C367").Value = "Coquimbo"
$ws.Range("D367").Value = 44610
$ws.Range("E367").Value = 5
$ws.Range("F367").Value = 100112003
$ws.Range("G367").Value = "Ajo"
$ws.Range("H367").Value = "Chino"
$ws.Range("I367").Value = "1a (cosecha)"
$ws.Range("J367").Value = 60
$ws.Range("K367").Value = 6000
$ws.Range("L367").Value = 6000
$ws.Range("M367").Value = 6000
$ws.Range("N367").Value = "`$/trenza 50 unidades"
$ws.Range("O367").Value = "Limache"
$ws.Range("P367").Value = 1200
$ws.Range("Q367").Value = 5
$ws.Range("R367").Value = "Hortaliza"

# Row 368: new "Limache" 2a (cosecha) entry, $/trenza 50 unidades
$ws.Range("A368").Value = 3
$ws.Range("B368").Value = "Femacal de La Calera"
$ws.Range("C368").Value = "Coquimbo"
$ws.Range("D368").Value = 44610
$ws.Range("E368").Value = 5
$ws.Range("F368").Value = 100112003
$ws.Range("G368").Value = "Ajo"
$ws.Range("H368").Value = "Chino"
$ws.Range("I368").Value = "2a (cosecha)"
$ws.Range("J368").Value = 55
$ws.Range("K368").Value = 4000
$ws.Range("L368").Value = 4000
$ws.Range("M368").Value = 4000
$ws.Range("N368").Value = "`$/trenza 50 unidades"
$ws.Range("O368").Value = "Limache"
$ws.Range("P368").Value = 800
$ws.Range("Q368").Value = 5
$ws.Range("R368").Value = "Hortaliza"

# Make sure date cells keep the date number format used by the rest of column D
$ws.Range("D366:D368").NumberFormat = $ws.Range("D365").NumberFormat
